# Actualización automática 2025-09-29 08:50:10
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M16").Value = 829.4400000000001
$ws1.Range("M26").Value = "6 de 24"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F16").Value = 829.4400000000001
$ws2.Range("F26").Value = 20529.03

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D12").Value = 20277.11
$ws3.Range("E12").Value = 22822.9754117774
$ws3.Range("F12").Value = 0.4704656569998151

$ws3.Range("D15").Value = 20529.03
$ws3.Range("E15").Value = 37674.43623249458
$ws3.Range("F15").Value = 0.3527114676984442
